$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor rows delivered 25 Jun 24

# Row 13: PT11-SPF
$ws.Range("A13").Value = "PT11-SPF"
$ws.Range("B13").Value = "6444fa9dbca6e305b2e2c466"
$ws.Range("C13").Value = "6444fa9e4426821a1c8b91a9"
$ws.Range("D13").Value = 45.349024
$ws.Range("E13").Value = -108.892649
$ws.Range("F13").NumberFormat = "d-mmm-yy"
$ws.Range("F13").Value = 45468

# Row 14: WS34-YKL
$ws.Range("A14").Value = "WS34-YKL"
$ws.Range("B14").Value = "6601c717158679306b7ee056"
$ws.Range("C14").Value = "6601c717f7bd28322857b8bb"
$ws.Range("D14").Value = 45.660301
$ws.Range("E14").Value = -108.850004
$ws.Range("F14").NumberFormat = "d-mmm-yy"
$ws.Range("F14").Value = 45468

# Row 15: WS35-F4R
$ws.Range("A15").Value = "WS35-F4R"
$ws.Range("B15").Value = "6601c966ca0eeb3251d33804"
$ws.Range("C15").Value = "6601c9664bbb713194d31a6a"
$ws.Range("D15").Value = 45.431862
$ws.Range("E15").Value = -108.882907
$ws.Range("F15").NumberFormat = "d-mmm-yy"
$ws.Range("F15").Value = 45468

# Row 16: WS36-DK3
$ws.Range("A16").Value = "WS36-DK3"
$ws.Range("B16").Value = "6601cbf05fe5e133cd1dd980"
$ws.Range("C16").Value = "6601cbf1022d5e321084e51a"
$ws.Range("D16").Value = 45.502678
$ws.Range("E16").Value = -108.861688
$ws.Range("F16").NumberFormat = "d-mmm-yy"
$ws.Range("F16").Value = 45468

# Update selection to reflect the saved view state
$ws.Range("B23").Select()
